$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header change
$ws.Range("A1").Value = "Estandar"

# Row 2
$ws.Range("B2").Value = 1001
$ws.Range("C2").Value = "2025-02-14 09:20:18"
$ws.Range("D2").Value = "A.VAZQUEZ"
$ws.Range("E2").Value = "2025-02-14 09:25:59"
$ws.Range("F2").Value = "A.VAZQUEZ"

# Row 3
$ws.Range("A3").Value = "ALCOHOL"
$ws.Range("B3").Value = 1002
$ws.Range("C3").Value = "2025-02-14 09:25:38"
$ws.Range("D3").Value = "N.SALINAS"
$ws.Range("E3").Value = "2025-02-14 09:30:21"
$ws.Range("F3").Value = "A.VAZQUEZ"

# Row 4
$ws.Range("A4").Value = "ACIDO"
$ws.Range("B4").Value = 1009
$ws.Range("C4").Value = "2025-02-14 09:26:28"
$ws.Range("D4").Value = "N.SALINAS"
$ws.Range("E4").Value = "2025-02-14 09:30:40"
$ws.Range("F4").Value = "A.VAZQUEZ"

# Row 5
$ws.Range("A5").Value = "ACIDO"
$ws.Range("B5").Value = 1009
$ws.Range("C5").Value = "2025-02-14 09:30:54"
$ws.Range("D5").Value = "A.VAZQUEZ"
$ws.Range("E5").Value = "2025-02-14 09:45:19"
$ws.Range("F5").Value = "A.VAZQUEZ"

# Row 6 (new)
$ws.Range("A6").Value = "ALCOHOL"
$ws.Range("B6").Value = 1006
$ws.Range("C6").Value = "2025-02-14 09:31:35"
$ws.Range("D6").Value = "N.SALINAS"
$ws.Range("E6").Value = "2025-02-14 09:46:00"
$ws.Range("F6").Value = "N.SALINAS"

# Row 7 (new)
$ws.Range("A7").Value = "METANOL"
$ws.Range("B7").Value = 1007
$ws.Range("C7").Value = "2025-02-14 09:45:44"
$ws.Range("D7").Value = "A.VAZQUEZ"

# Row 8 (new)
$ws.Range("A8").Value = "ACETONITRILO"
$ws.Range("B8").Value = 1004
$ws.Range("C8").Value = "2025-02-14 09:46:16"
$ws.Range("D8").Value = "A.VAZQUEZ"

# Column width adjustments: split col A/B into separate widths
# (Column A keeps ~12.29 chars, column B widens to ~13.43 chars)
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666

# Selection update
$ws.Range("F8").Select()

$wb.Save()
